$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (existing) - update text, clear A1's style (it loses its style index)
$ws.Range("A1").Value = "Meta Serevice 1"
$ws.Range("A1").ClearFormats()
$ws.Range("B1").Value = 12.3

# Row 2
$ws.Range("A2").Value = "Meta Service 2"
$ws.Range("B2").Value = 9.2

# Row 3
$ws.Range("A3").Value = "AWS Service 1"
$ws.Range("B3").Value = 8.43

# Row 4
$ws.Range("A4").Value = "AWS Service 2"
$ws.Range("B4").Value = 9.43

# Apply same style as B1 (style index 1) to the rest of column B
$ws.Range("B2:B4").NumberFormat = "General"

# Column B width
$ws.Columns.Item(2).ColumnWidth = 10.83203125

# Selection / active cell
$ws.Range("G6").Select()
